# Refresh the "cryptos" price/volume snapshot to the latest scrape.
# Rows 21/22 and 36/37 also swap rank order (Dai<->Avalanche, ImmutableX<->RenderToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.158.46"
$ws.Range("E2").Value = "  +1.18%  "
# Row 3
$ws.Range("D3").Value = "1.787.90"
$ws.Range("E3").Value = "  +1.27%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").Value = "'226.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
# Row 6
$ws.Range("D6").Value = "'0.547"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "
# Row 7
$ws.Range("E7").Value = "  -0.04%  "
# Row 8
$ws.Range("D8").Value = "'31.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.97%  "
# Row 9
$ws.Range("D9").Value = "'0.291"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "
# Row 10
$ws.Range("D10").Value = "'0.0689"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.22%  "
# Row 11
$ws.Range("D11").Value = "'0.0947"
$ws.Range("D11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "2.046.51"
$ws.Range("E12").Value = "  +1.23%  "
# Row 13
$ws.Range("D13").Value = "'11.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.33%  "
# Row 14
$ws.Range("D14").Value = "1.785.54"
$ws.Range("E14").Value = "  +1.39%  "
# Row 15
$ws.Range("D15").Value = "'0.623"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.82%  "
# Row 16
$ws.Range("D16").Value = "34.114.97"
$ws.Range("E16").Value = "  +1.16%  "
# Row 17
$ws.Range("E17").Value = "  +1.22%  "
# Row 18
$ws.Range("D18").Value = "'68.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.43%  "
# Row 19
$ws.Range("D19").Value = "'247.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.98%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0777"
$ws.Range("E20").Value = "  +0.51%  "
# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'10.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.21%  "
# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
# Row 23
$ws.Range("E23").Value = "  +1.46%  "
# Row 24
$ws.Range("D24").Value = "'2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "
# Row 25
$ws.Range("D25").Value = "'161.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
# Row 26
$ws.Range("E26").Value = "  +2.16%  "
# Row 27
$ws.Range("E27").Value = "  +1.47%  "
# Row 28
$ws.Range("E28").Value = "  +1.31%  "
# Row 29
$ws.Range("E29").Value = "  +0.10%  "
# Row 30
$ws.Range("E30").Value = "  +0.20%  "
# Row 31
$ws.Range("E31").Value = "  +2.19%  "
# Row 32
$ws.Range("D32").Value = "'3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.23%  "
# Row 33
$ws.Range("D33").Value = "'3.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.98%  "
# Row 34
$ws.Range("E34").Value = "  +1.24%  "
# Row 35
$ws.Range("D35").Value = "1.447.34"
$ws.Range("E35").Value = "  +4.93%  "
# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.30%  "
# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.654"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
# Row 38
$ws.Range("D38").Value = "'0.0191"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.49%  "
# Row 39
$ws.Range("E39").Value = "  +0.75%  "
# Row 40
$ws.Range("E40").Value = "  +3.91%  "
# Row 41
$ws.Range("E41").Value = "  +0.77%  "
# Row 42
$ws.Range("D42").Value = "'0.922"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
# Row 43
$ws.Range("E43").Value = "  +0.98%  "
# Row 44
$ws.Range("D44").Value = "'13.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.68%  "
# Row 45
$ws.Range("E45").Value = "  +4.37%  "
# Row 46
$ws.Range("E46").Value = "  +2.10%  "
# Row 47
$ws.Range("D47").Value = "'1.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
# Row 48
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  -0.69%  "
# Row 49
$ws.Range("D49").Value = "1.947.91"
$ws.Range("E49").Value = "  +1.57%  "
# Row 50
$ws.Range("D50").Value = "'105.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "
# Row 51
$ws.Range("E51").Value = "  -0.08%  "
